$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.734.11"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "3.036.13"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0862"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "3.526.48"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.025.30"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.979"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -13.41%  "
$ws.Range("D19").Value = "51.737.54"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("E27").Value = "  +6.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.175"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0450"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.53%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.116"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +7.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").Value = "2.038.10"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "3.336.84"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.25%  "
